$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "Wins", "Losses", "Ties" columns, copying the
# existing header formatting (bold, bordered, centered) from column AC.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the season record (Wins / Losses / Ties) for every player row.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 103
    $ws.Cells.Item($row, 31).Value = 59
    $ws.Cells.Item($row, 32).Value = 0
}
